# "Handles float input without breaking stuff"
#
# Re-grades the marksheet: updates the Right/Wrong/Not-Attempt/Max summary
# block (rows 10-12), fixes the "Marking" wrong-answer penalty to be a real
# number instead of text, records the final score as a fraction string, and
# fills in the actual per-question "Student Ans" column (A16:A40) with the
# student's picked option -- colored green/red/plain depending on whether it
# matches the "Correct Ans" column (B16:B40). The two extra (unused) answer
# blocks in columns D:E and G:H are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----------------------------------------
# No. (attempt counts)
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

# Marking (per-question weight) - C11 was stored as text "-1"; make it numeric
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Total
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "60/112"

# The row-labels (No./Marking/Total) pick up the same centered/bold style
# used elsewhere in the sheet (mtitleStyle, already applied to A9/A15).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Student answers (column A, rows 16-40) ----------------------------
# Style index 5 ("correctStyle"/green) when the student's answer matches the
# correct answer, style index 6 ("incorrectStyle"/red) when it doesn't, and
# the plain style 7 (already in place) when left blank / not attempted.
# Pull the format from cells that already use those exact styles (B10 = 5,
# C10 = 6) so no new style entries get minted.

function Set-StudentAnswer($row, $value, $correct) {
    $cell = $ws.Range("A" + $row)
    if ($value -eq $correct) {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $cell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $cell.Value = $value
}

Set-StudentAnswer 16 "Option A" "Option A"
Set-StudentAnswer 18 "Option B" "Option B"
Set-StudentAnswer 19 "Option B" "Option C"
Set-StudentAnswer 22 "Option D" "Option D"
Set-StudentAnswer 24 "Option A" "Option A"
Set-StudentAnswer 25 "Option A" "Option A"
Set-StudentAnswer 26 "Option D" "Option C"
Set-StudentAnswer 27 "Option A" "Option A"
Set-StudentAnswer 28 "Option D" "Option D"
Set-StudentAnswer 30 "Option B" "Option B"
Set-StudentAnswer 31 "Option D" "Option D"
Set-StudentAnswer 32 "Option C" "Option C"
Set-StudentAnswer 33 "Option A" "Option D"
Set-StudentAnswer 34 "Option C" "Option B"
Set-StudentAnswer 35 "Option D" "Option D"
Set-StudentAnswer 37 "Option A" "Option A"
Set-StudentAnswer 38 "Option A" "Option A"
Set-StudentAnswer 39 "Option D" "Option D"

# Rows 17, 20, 21, 23, 29, 36, 40 stay blank/not-attempted (already style 7).

# ---- Drop the two unused answer blocks ---------------------------------
$ws.Range("D16:E40").Clear()
$ws.Range("G15:H21").Clear()
